$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 36
$ws.Range("E2").Value = 28
$ws.Range("G2").Value = 0.778
$ws.Range("I2").Value = 37.6
$ws.Range("J2").Value = 80.5
$ws.Range("L2").Value = 9.300000000000001
$ws.Range("M2").Value = 24.7
$ws.Range("N2").Value = 0.378
$ws.Range("P2").Value = 23.3
$ws.Range("Q2").Value = 0.77
$ws.Range("R2").Value = 8.4
$ws.Range("S2").Value = 32.9
$ws.Range("U2").Value = 25.2
$ws.Range("V2").Value = 14.6
$ws.Range("W2").Value = 8.800000000000001
$ws.Range("X2").Value = 4.7
$ws.Range("Y2").Value = 4.8
$ws.Range("Z2").Value = 18.4
$ws.Range("AA2").Value = 21.3
$ws.Range("AB2").Value = 102.5
$ws.Range("AC2").Value = 5.1
$ws.Range("AD2").Value = 20
$ws.Range("AE2").Value = 3
$ws.Range("AH2").Value = 21
$ws.Range("AN2").Value = 4
$ws.Range("AS2").Value = 12
$ws.Range("AU2").Value = 3
$ws.Range("AV2").Value = 15
$ws.Range("AX2").Value = 15
$ws.Range("AY2").Value = 15
$ws.Range("BF2").Value = "2015-01-11"
$ws.Range("AD3").Value = 29
$ws.Range("AH3").Value = 4
$ws.Range("AR3").Value = 20
$ws.Range("AT3").Value = 12
$ws.Range("AU3").Value = 2
$ws.Range("AV3").Value = 17
$ws.Range("AZ3").Value = 20
$ws.Range("BA3").Value = 29
$ws.Range("BF3").Value = "2015-01-11"
$ws.Range("AD4").Value = 11
$ws.Range("AI4").Value = 24
$ws.Range("AS4").Value = 13
$ws.Range("AT4").Value = 16
$ws.Range("AU4").Value = 19
$ws.Range("AV4").Value = 18
$ws.Range("AW4").Value = 22
$ws.Range("AX4").Value = 22
$ws.Range("BF4").Value = "2015-01-11"
$ws.Range("AE5").Value = 19
$ws.Range("AJ5").Value = 12
$ws.Range("AO5").Value = 17
$ws.Range("BA5").Value = 8
$ws.Range("BF5").Value = "2015-01-11"
$ws.Range("AH6").Value = 11
$ws.Range("AI6").Value = 21
$ws.Range("AN6").Value = 13
$ws.Range("BB6").Value = 12
$ws.Range("BF6").Value = "2015-01-11"
$ws.Range("D7").Value = 37
$ws.Range("F7").Value = 18
$ws.Range("G7").Value = 0.514
$ws.Range("I7").Value = 36.7
$ws.Range("J7").Value = 81.7
$ws.Range("K7").Value = 0.449
$ws.Range("M7").Value = 23.4
$ws.Range("O7").Value = 18.8
$ws.Range("P7").Value = 24.7
$ws.Range("Q7").Value = 0.762
$ws.Range("R7").Value = 10.9
$ws.Range("S7").Value = 30.2
$ws.Range("U7").Value = 22.1
$ws.Range("W7").Value = 7.2
$ws.Range("AA7").Value = 21.5
$ws.Range("AB7").Value = 100.2
$ws.Range("AC7").Value = 0.4
$ws.Range("AD7").Value = 11
$ws.Range("AF7").Value = 13
$ws.Range("AG7").Value = 13
$ws.Range("AI7").Value = 18
$ws.Range("AL7").Value = 11
$ws.Range("AM7").Value = 12
$ws.Range("AN7").Value = 20
$ws.Range("AP7").Value = 11
$ws.Range("AQ7").Value = 11
$ws.Range("AR7").Value = 15
$ws.Range("AT7").Value = 25
$ws.Range("AW7").Value = 20
$ws.Range("AX7").Value = 27
$ws.Range("AY7").Value = 19
$ws.Range("BA7").Value = 7
$ws.Range("BF7").Value = "2015-01-11"
$ws.Range("D8").Value = 38
$ws.Range("F8").Value = 12
$ws.Range("G8").Value = 0.6840000000000001
$ws.Range("I8").Value = 41
$ws.Range("J8").Value = 86.3
$ws.Range("L8").Value = 9.6
$ws.Range("M8").Value = 26.7
$ws.Range("P8").Value = 22.6
$ws.Range("Q8").Value = 0.762
$ws.Range("R8").Value = 10.8
$ws.Range("S8").Value = 31.4
$ws.Range("T8").Value = 42.2
$ws.Range("U8").Value = 23.7
$ws.Range("W8").Value = 8
$ws.Range("AB8").Value = 108.8
$ws.Range("AC8").Value = 6.3
$ws.Range("AD8").Value = 7
$ws.Range("AF8").Value = 8
$ws.Range("AG8").Value = 8
$ws.Range("AH8").Value = 11
$ws.Range("AJ8").Value = 5
$ws.Range("AM8").Value = 5
$ws.Range("AO8").Value = 16
$ws.Range("AP8").Value = 18
$ws.Range("AR8").Value = 16
$ws.Range("AT8").Value = 19
$ws.Range("BA8").Value = 4
$ws.Range("BC8").Value = 4
$ws.Range("BF8").Value = "2015-01-11"
$ws.Range("AD9").Value = 11
$ws.Range("AE9").Value = 17
$ws.Range("AH9").Value = 15
$ws.Range("AJ9").Value = 4
$ws.Range("AK9").Value = 23
$ws.Range("AM9").Value = 11
$ws.Range("AV9").Value = 16
$ws.Range("BB9").Value = 9
$ws.Range("BF9").Value = "2015-01-11"
$ws.Range("AD10").Value = 11
$ws.Range("AH10").Value = 15
$ws.Range("AJ10").Value = 6
$ws.Range("AN10").Value = 18
$ws.Range("AY10").Value = 16
$ws.Range("BF10").Value = "2015-01-11"
$ws.Range("AE11").Value = 1
$ws.Range("AN11").Value = 3
$ws.Range("BA11").Value = 27
$ws.Range("BF11").Value = "2015-01-11"
$ws.Range("D12").Value = 37
$ws.Range("E12").Value = 26
$ws.Range("G12").Value = 0.703
$ws.Range("H12").Value = 48.5
$ws.Range("I12").Value = 36.1
$ws.Range("J12").Value = 83
$ws.Range("M12").Value = 33.8
$ws.Range("N12").Value = 0.349
$ws.Range("O12").Value = 17.6
$ws.Range("P12").Value = 24.9
$ws.Range("Q12").Value = 0.706
$ws.Range("S12").Value = 31.5
$ws.Range("T12").Value = 43.9
$ws.Range("U12").Value = 20.7
$ws.Range("V12").Value = 17.4
$ws.Range("W12").Value = 9.699999999999999
$ws.Range("Y12").Value = 5.6
$ws.Range("Z12").Value = 23.1
$ws.Range("AA12").Value = 20.9
$ws.Range("AB12").Value = 101.5
$ws.Range("AC12").Value = 4.7
$ws.Range("AD12").Value = 11
$ws.Range("AE12").Value = 4
$ws.Range("AG12").Value = 4
$ws.Range("AI12").Value = 25
$ws.Range("AK12").Value = 26
$ws.Range("AO12").Value = 14
$ws.Range("AP12").Value = 8
$ws.Range("AS12").Value = 20
$ws.Range("AT12").Value = 9
$ws.Range("AU12").Value = 22
$ws.Range("AX12").Value = 16
$ws.Range("AY12").Value = 25
$ws.Range("AZ12").Value = 29
$ws.Range("BA12").Value = 14
$ws.Range("BB12").Value = 14
$ws.Range("BF12").Value = "2015-01-11"
$ws.Range("AE13").Value = 19
$ws.Range("AL13").Value = 19
$ws.Range("AR13").Value = 13
$ws.Range("AU13").Value = 21
$ws.Range("AX13").Value = 18
$ws.Range("AY13").Value = 18
$ws.Range("BA13").Value = 11
$ws.Range("BF13").Value = "2015-01-11"
$ws.Range("D14").Value = 37
$ws.Range("F14").Value = 12
$ws.Range("G14").Value = 0.676
$ws.Range("I14").Value = 38.8
$ws.Range("J14").Value = 81.90000000000001
$ws.Range("K14").Value = 0.474
$ws.Range("N14").Value = 0.385
$ws.Range("O14").Value = 18.9
$ws.Range("P14").Value = 25.2
$ws.Range("R14").Value = 8.699999999999999
$ws.Range("S14").Value = 32.6
$ws.Range("T14").Value = 41.2
$ws.Range("U14").Value = 24.9
$ws.Range("W14").Value = 7.8
$ws.Range("Z14").Value = 20.5
$ws.Range("AA14").Value = 21.9
$ws.Range("AB14").Value = 106.8
$ws.Range("AC14").Value = 7
$ws.Range("AD14").Value = 11
$ws.Range("AF14").Value = 8
$ws.Range("AI14").Value = 9
$ws.Range("AJ14").Value = 20
$ws.Range("AM14").Value = 4
$ws.Range("AO14").Value = 5
$ws.Range("AS14").Value = 14
$ws.Range("AT14").Value = 24
$ws.Range("AV14").Value = 4
$ws.Range("AX14").Value = 13
$ws.Range("BC14").Value = 3
$ws.Range("BF14").Value = "2015-01-11"
$ws.Range("I15").Value = 37.9
$ws.Range("J15").Value = 86.5
$ws.Range("K15").Value = 0.438
$ws.Range("L15").Value = 7.1
$ws.Range("M15").Value = 19.8
$ws.Range("O15").Value = 18.8
$ws.Range("P15").Value = 25.2
$ws.Range("Q15").Value = 0.747
$ws.Range("R15").Value = 11.7
$ws.Range("S15").Value = 31.3
$ws.Range("T15").Value = 43.1
$ws.Range("U15").Value = 20.6
$ws.Range("W15").Value = 7.5
$ws.Range("Z15").Value = 21.8
$ws.Range("AA15").Value = 20.7
$ws.Range("AB15").Value = 101.6
$ws.Range("AC15").Value = -5.9
$ws.Range("AD15").Value = 11
$ws.Range("AH15").Value = 15
$ws.Range("AK15").Value = 24
$ws.Range("AL15").Value = 20
$ws.Range("AN15").Value = 12
$ws.Range("AP15").Value = 7
$ws.Range("AR15").Value = 7
$ws.Range("AS15").Value = 22
$ws.Range("AT15").Value = 14
$ws.Range("AU15").Value = 23
$ws.Range("AW15").Value = 17
$ws.Range("AZ15").Value = 22
$ws.Range("BA15").Value = 16
$ws.Range("BB15").Value = 13
$ws.Range("BF15").Value = "2015-01-11"
$ws.Range("D16").Value = 36
$ws.Range("E16").Value = 25
$ws.Range("G16").Value = 0.694
$ws.Range("H16").Value = 49
$ws.Range("I16").Value = 38.8
$ws.Range("J16").Value = 83.5
$ws.Range("K16").Value = 0.464
$ws.Range("L16").Value = 5.6
$ws.Range("M16").Value = 15.9
$ws.Range("N16").Value = 0.356
$ws.Range("Q16").Value = 0.776
$ws.Range("S16").Value = 31.7
$ws.Range("T16").Value = 42.1
$ws.Range("V16").Value = 12.8
$ws.Range("X16").Value = 4.3
$ws.Range("Z16").Value = 19.4
$ws.Range("AA16").Value = 20.5
$ws.Range("AB16").Value = 101.2
$ws.Range("AC16").Value = 3.5
$ws.Range("AD16").Value = 20
$ws.Range("AE16").Value = 7
$ws.Range("AG16").Value = 5
$ws.Range("AH16").Value = 2
$ws.Range("AI16").Value = 10
$ws.Range("AN16").Value = 14
$ws.Range("AO16").Value = 9
$ws.Range("AQ16").Value = 7
$ws.Range("AS16").Value = 18
$ws.Range("AT16").Value = 20
$ws.Range("AX16").Value = 24
$ws.Range("AY16").Value = 23
$ws.Range("AZ16").Value = 8
$ws.Range("BA16").Value = 20
$ws.Range("BB16").Value = 16
$ws.Range("BF16").Value = "2015-01-11"
$ws.Range("D17").Value = 36
$ws.Range("E17").Value = 15
$ws.Range("G17").Value = 0.417
$ws.Range("I17").Value = 34.2
$ws.Range("J17").Value = 74.09999999999999
$ws.Range("K17").Value = 0.461
$ws.Range("L17").Value = 7.4
$ws.Range("N17").Value = 0.352
$ws.Range("Q17").Value = 0.73
$ws.Range("R17").Value = 8.199999999999999
$ws.Range("S17").Value = 28
$ws.Range("T17").Value = 36.2
$ws.Range("U17").Value = 19.7
$ws.Range("Y17").Value = 4.4
$ws.Range("AB17").Value = 93.8
$ws.Range("AC17").Value = -4.5
$ws.Range("AD17").Value = 20
$ws.Range("AE17").Value = 19
$ws.Range("AG17").Value = 19
$ws.Range("AO17").Value = 10
$ws.Range("AP17").Value = 9
$ws.Range("AU17").Value = 30
$ws.Range("BA17").Value = 15
$ws.Range("BF17").Value = "2015-01-11"
$ws.Range("AF18").Value = 15
$ws.Range("AG18").Value = 14
$ws.Range("AH18").Value = 8
$ws.Range("AT18").Value = 26
$ws.Range("AZ18").Value = 27
$ws.Range("BF18").Value = "2015-01-11"
$ws.Range("AD19").Value = 20
$ws.Range("AH19").Value = 21
$ws.Range("AK19").Value = 25
$ws.Range("AO19").Value = 6
$ws.Range("AQ19").Value = 24
$ws.Range("AT19").Value = 27
$ws.Range("AX19").Value = 23
$ws.Range("BF19").Value = "2015-01-11"
$ws.Range("D20").Value = 36
$ws.Range("E20").Value = 18
$ws.Range("G20").Value = 0.5
$ws.Range("J20").Value = 84.8
$ws.Range("K20").Value = 0.458
$ws.Range("L20").Value = 6.6
$ws.Range("M20").Value = 19.4
$ws.Range("N20").Value = 0.342
$ws.Range("O20").Value = 17
$ws.Range("P20").Value = 22.5
$ws.Range("Q20").Value = 0.755
$ws.Range("R20").Value = 11.7
$ws.Range("S20").Value = 31.8
$ws.Range("T20").Value = 43.5
$ws.Range("V20").Value = 12.7
$ws.Range("X20").Value = 6
$ws.Range("Y20").Value = 5.8
$ws.Range("AA20").Value = 18.9
$ws.Range("AC20").Value = 0.5
$ws.Range("AD20").Value = 20
$ws.Range("AE20").Value = 15
$ws.Range("AF20").Value = 13
$ws.Range("AG20").Value = 15
$ws.Range("AI20").Value = 8
$ws.Range("AJ20").Value = 11
$ws.Range("AK20").Value = 13
$ws.Range("AN20").Value = 21
$ws.Range("AP20").Value = 19
$ws.Range("AQ20").Value = 14
$ws.Range("AR20").Value = 8
$ws.Range("AS20").Value = 17
$ws.Range("AT20").Value = 11
$ws.Range("AW20").Value = 18
$ws.Range("AX20").Value = 3
$ws.Range("AY20").Value = 27
$ws.Range("AZ20").Value = 6
$ws.Range("BA20").Value = 26
$ws.Range("BB20").Value = 15
$ws.Range("BF20").Value = "2015-01-11"
$ws.Range("AJ21").Value = 21
$ws.Range("AQ21").Value = 8
$ws.Range("BA21").Value = 28
$ws.Range("BF21").Value = "2015-01-11"
$ws.Range("AD22").Value = 11
$ws.Range("AF22").Value = 15
$ws.Range("AG22").Value = 16
$ws.Range("AH22").Value = 15
$ws.Range("AK22").Value = 22
$ws.Range("AU22").Value = 28
$ws.Range("AV22").Value = 25
$ws.Range("AX22").Value = 4
$ws.Range("AZ22").Value = 28
$ws.Range("BA22").Value = 19
$ws.Range("BF22").Value = "2015-01-11"
$ws.Range("D23").Value = 40
$ws.Range("F23").Value = 27
$ws.Range("G23").Value = 0.325
$ws.Range("I23").Value = 36.5
$ws.Range("J23").Value = 81.09999999999999
$ws.Range("K23").Value = 0.45
$ws.Range("L23").Value = 7
$ws.Range("M23").Value = 19
$ws.Range("O23").Value = 13.7
$ws.Range("P23").Value = 18.8
$ws.Range("Q23").Value = 0.728
$ws.Range("R23").Value = 9
$ws.Range("S23").Value = 32
$ws.Range("T23").Value = 41
$ws.Range("U23").Value = 19.9
$ws.Range("V23").Value = 14.7
$ws.Range("W23").Value = 7.1
$ws.Range("X23").Value = 3.6
$ws.Range("Y23").Value = 5.8
$ws.Range("Z23").Value = 21.2
$ws.Range("AA23").Value = 18.3
$ws.Range("AB23").Value = 93.59999999999999
$ws.Range("AC23").Value = -5.8
$ws.Range("AD23").Value = 1
$ws.Range("AI23").Value = 22
$ws.Range("AK23").Value = 17
$ws.Range("AM23").Value = 25
$ws.Range("AS23").Value = 16
$ws.Range("AU23").Value = 27
$ws.Range("AV23").Value = 19
$ws.Range("AW23").Value = 21
$ws.Range("AY23").Value = 26
$ws.Range("BA23").Value = 30
$ws.Range("BF23").Value = "2015-01-11"
$ws.Range("AD24").Value = 20
$ws.Range("AH24").Value = 21
$ws.Range("AL24").Value = 18
$ws.Range("AP24").Value = 9
$ws.Range("AR24").Value = 9
$ws.Range("AU24").Value = 26
$ws.Range("AX24").Value = 6
$ws.Range("AY24").Value = 24
$ws.Range("BF24").Value = "2015-01-11"
$ws.Range("D25").Value = 39
$ws.Range("F25").Value = 17
$ws.Range("G25").Value = 0.5639999999999999
$ws.Range("H25").Value = 48.5
$ws.Range("I25").Value = 39.8
$ws.Range("J25").Value = 85.90000000000001
$ws.Range("K25").Value = 0.463
$ws.Range("M25").Value = 27.1
$ws.Range("N25").Value = 0.366
$ws.Range("Q25").Value = 0.798
$ws.Range("R25").Value = 10.7
$ws.Range("S25").Value = 31.7
$ws.Range("T25").Value = 42.4
$ws.Range("V25").Value = 15.2
$ws.Range("W25").Value = 8.6
$ws.Range("X25").Value = 4.9
$ws.Range("Y25").Value = 3.9
$ws.Range("Z25").Value = 22.6
$ws.Range("AA25").Value = 21
$ws.Range("AB25").Value = 106.9
$ws.Range("AC25").Value = 2.3
$ws.Range("AD25").Value = 3
$ws.Range("AE25").Value = 12
$ws.Range("AH25").Value = 13
$ws.Range("AJ25").Value = 7
$ws.Range("AN25").Value = 8
$ws.Range("AO25").Value = 15
$ws.Range("AS25").Value = 19
$ws.Range("AT25").Value = 17
$ws.Range("AU25").Value = 20
$ws.Range("AV25").Value = 24
$ws.Range("AZ25").Value = 26
$ws.Range("BA25").Value = 12
$ws.Range("BC25").Value = 12
$ws.Range("BF25").Value = "2015-01-11"
$ws.Range("D26").Value = 37
$ws.Range("E26").Value = 29
$ws.Range("G26").Value = 0.784
$ws.Range("M26").Value = 27
$ws.Range("N26").Value = 0.376
$ws.Range("O26").Value = 15.4
$ws.Range("P26").Value = 19.5
$ws.Range("Q26").Value = 0.792
$ws.Range("R26").Value = 11.3
$ws.Range("T26").Value = 46.6
$ws.Range("U26").Value = 22.8
$ws.Range("V26").Value = 14.1
$ws.Range("X26").Value = 5.5
$ws.Range("AA26").Value = 19.1
$ws.Range("AB26").Value = 103.7
$ws.Range("AC26").Value = 7.3
$ws.Range("AD26").Value = 11
$ws.Range("AH26").Value = 6
$ws.Range("AK26").Value = 16
$ws.Range("AR26").Value = 12
$ws.Range("AX26").Value = 7
$ws.Range("AZ26").Value = 7
$ws.Range("BF26").Value = "2015-01-11"
$ws.Range("D27").Value = 36
$ws.Range("E27").Value = 15
$ws.Range("G27").Value = 0.417
$ws.Range("I27").Value = 36.6
$ws.Range("J27").Value = 79.8
$ws.Range("K27").Value = 0.458
$ws.Range("L27").Value = 5.3
$ws.Range("O27").Value = 23.7
$ws.Range("P27").Value = 30.6
$ws.Range("R27").Value = 11.1
$ws.Range("S27").Value = 33.4
$ws.Range("V27").Value = 16.5
$ws.Range("W27").Value = 6.4
$ws.Range("Y27").Value = 6.3
$ws.Range("Z27").Value = 22.1
$ws.Range("AA27").Value = 25.6
$ws.Range("AC27").Value = -2
$ws.Range("AD27").Value = 20
$ws.Range("AE27").Value = 19
$ws.Range("AG27").Value = 19
$ws.Range("AH27").Value = 5
$ws.Range("AI27").Value = 20
$ws.Range("AK27").Value = 12
$ws.Range("AM27").Value = 28
$ws.Range("AN27").Value = 19
$ws.Range("AX27").Value = 26
$ws.Range("BF27").Value = "2015-01-11"
$ws.Range("D28").Value = 38
$ws.Range("E28").Value = 23
$ws.Range("G28").Value = 0.605
$ws.Range("J28").Value = 82.5
$ws.Range("K28").Value = 0.467
$ws.Range("L28").Value = 8.199999999999999
$ws.Range("M28").Value = 21.9
$ws.Range("N28").Value = 0.375
$ws.Range("O28").Value = 17
$ws.Range("P28").Value = 22.4
$ws.Range("Q28").Value = 0.758
$ws.Range("S28").Value = 34
$ws.Range("T28").Value = 43.8
$ws.Range("V28").Value = 14.8
$ws.Range("Y28").Value = 4.4
$ws.Range("Z28").Value = 19.6
$ws.Range("AD28").Value = 7
$ws.Range("AJ28").Value = 19
$ws.Range("AL28").Value = 10
$ws.Range("AO28").Value = 18
$ws.Range("AP28").Value = 20
$ws.Range("AQ28").Value = 13
$ws.Range("AT28").Value = 10
$ws.Range("AV28").Value = 20
$ws.Range("AX28").Value = 8
$ws.Range("AZ28").Value = 9
$ws.Range("BA28").Value = 18
$ws.Range("BB28").Value = 10
$ws.Range("BF28").Value = "2015-01-11"
$ws.Range("AD29").Value = 20
$ws.Range("AG29").Value = 5
$ws.Range("AH29").Value = 9
$ws.Range("AK29").Value = 11
$ws.Range("AN29").Value = 9
$ws.Range("AR29").Value = 11
$ws.Range("AW29").Value = 19
$ws.Range("AZ29").Value = 21
$ws.Range("BA29").Value = 3
$ws.Range("BF29").Value = "2015-01-11"
$ws.Range("AR30").Value = 10
$ws.Range("AT30").Value = 18
$ws.Range("AW30").Value = 26
$ws.Range("AX30").Value = 5
$ws.Range("BF30").Value = "2015-01-11"
$ws.Range("D31").Value = 36
$ws.Range("F31").Value = 11
$ws.Range("G31").Value = 0.694
$ws.Range("I31").Value = 39
$ws.Range("J31").Value = 82.5
$ws.Range("M31").Value = 15.5
$ws.Range("N31").Value = 0.397
$ws.Range("O31").Value = 15.7
$ws.Range("P31").Value = 21.3
$ws.Range("Q31").Value = 0.735
$ws.Range("R31").Value = 10.5
$ws.Range("T31").Value = 42.9
$ws.Range("V31").Value = 14.4
$ws.Range("Z31").Value = 21.3
$ws.Range("AA31").Value = 20.5
$ws.Range("AB31").Value = 99.90000000000001
$ws.Range("AC31").Value = 2.9
$ws.Range("AD31").Value = 20
$ws.Range("AF31").Value = 4
$ws.Range("AG31").Value = 5
$ws.Range("AI31").Value = 7
$ws.Range("AJ31").Value = 17
$ws.Range("AM31").Value = 29
$ws.Range("AQ31").Value = 25
$ws.Range("AR31").Value = 19
$ws.Range("AS31").Value = 15
$ws.Range("AV31").Value = 13
$ws.Range("AX31").Value = 10
$ws.Range("AY31").Value = 8
$ws.Range("BA31").Value = 17
$ws.Range("BC31").Value = 11
$ws.Range("BF31").Value = "2015-01-11"
